$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CategoryMap")

# Update individual data values (G column and D15) per diff
$ws.Range("G2").Value = 2
$ws.Range("G7").Value = 8
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 4
$ws.Range("G11").Value = 6
$ws.Range("G13").Value = 6
$ws.Range("D15").Value = 6

# Add a new blank row 20 (B:J), matching the existing body style (s="1", wrap text)
$ws.Range("B20:J20").WrapText = $true

# Update selection to match diff (select entire row 7)
$ws.Range("A7:XFD7").Select()

$wb.Save()
